$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1:E1").Value = 15
